# Generate Report for handback
#
# The localization-status report is being refreshed to reflect a completed
# handback: the Status text changes from "Ready for handoff" to
# "Handed back: in sync with en-US", and each localized-file row gains its
# "Latest Target File" / "Latest Handback File" columns (E/F) populated
# (hyperlinked, same as the existing A/C columns), plus an updated
# "Latest Handback DateTime" (column G).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Status column (B) for the two file rows - updates the shared Status text
$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

# Row 2 (2f835eaa-...md): populate Latest Target File / Latest Handback File
$ws.Hyperlinks.Add(
    $ws.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8b2153e47aa133a5f1540b3dc043b6d4ca42316c/e2e/2f835eaa-2126-4b04-835e-a4a65178257a.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "2f835eaa-2126-4b04-835e-a4a65178257a.md"
) | Out-Null
$ws.Range("E2").Font.Underline = $true
$ws.Range("E2").Font.Color = 15570276

$ws.Hyperlinks.Add(
    $ws.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34805292376e9622dc838562a4bd583e1922324b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/2f835eaa-2126-4b04-835e-a4a65178257a.639b1ad823e1f8a6516433323a7207e1bd38ba8c.zh-cn.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "2f835eaa-2126-4b04-835e-a4a65178257a.639b1ad823e1f8a6516433323a7207e1bd38ba8c.zh-cn.xlf"
) | Out-Null
$ws.Range("F2").Font.Underline = $true
$ws.Range("F2").Font.Color = 15570276

# Latest Handback DateTime for row 2
$ws.Range("G2").Value = "2016-01-18 07:08:52"

# Row 3 (845f499b-...md): populate Latest Target File / Latest Handback File
$ws.Hyperlinks.Add(
    $ws.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8b2153e47aa133a5f1540b3dc043b6d4ca42316c/e2e/845f499b-b28a-4a9b-aa7e-18df09bf383c.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "845f499b-b28a-4a9b-aa7e-18df09bf383c.md"
) | Out-Null
$ws.Range("E3").Font.Underline = $true
$ws.Range("E3").Font.Color = 15570276

$ws.Hyperlinks.Add(
    $ws.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34805292376e9622dc838562a4bd583e1922324b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/845f499b-b28a-4a9b-aa7e-18df09bf383c.cf05ccea05cfdc281aff4fe5748e88d52fbec87d.zh-cn.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "845f499b-b28a-4a9b-aa7e-18df09bf383c.cf05ccea05cfdc281aff4fe5748e88d52fbec87d.zh-cn.xlf"
) | Out-Null
$ws.Range("F3").Font.Underline = $true
$ws.Range("F3").Font.Color = 15570276

# Latest Handback DateTime for row 3
$ws.Range("G3").Value = "2016-01-18 07:08:52"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("B2").Value = $newStatus
$ws2.Range("B3").Value = $newStatus

$ws2.Hyperlinks.Add(
    $ws2.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8b2153e47aa133a5f1540b3dc043b6d4ca42316c/e2e/2f835eaa-2126-4b04-835e-a4a65178257a.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "2f835eaa-2126-4b04-835e-a4a65178257a.md"
) | Out-Null
$ws2.Range("E2").Font.Underline = $true
$ws2.Range("E2").Font.Color = 15570276

$ws2.Hyperlinks.Add(
    $ws2.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fccd6562d0a1ec584cbd387cc718777226d81b2c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/2f835eaa-2126-4b04-835e-a4a65178257a.639b1ad823e1f8a6516433323a7207e1bd38ba8c.de-de.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "2f835eaa-2126-4b04-835e-a4a65178257a.639b1ad823e1f8a6516433323a7207e1bd38ba8c.de-de.xlf"
) | Out-Null
$ws2.Range("F2").Font.Underline = $true
$ws2.Range("F2").Font.Color = 15570276

$ws2.Range("G2").Value = "2016-01-18 07:09:09"

$ws2.Hyperlinks.Add(
    $ws2.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8b2153e47aa133a5f1540b3dc043b6d4ca42316c/e2e/845f499b-b28a-4a9b-aa7e-18df09bf383c.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "845f499b-b28a-4a9b-aa7e-18df09bf383c.md"
) | Out-Null
$ws2.Range("E3").Font.Underline = $true
$ws2.Range("E3").Font.Color = 15570276

$ws2.Hyperlinks.Add(
    $ws2.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fccd6562d0a1ec584cbd387cc718777226d81b2c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/845f499b-b28a-4a9b-aa7e-18df09bf383c.cf05ccea05cfdc281aff4fe5748e88d52fbec87d.de-de.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "845f499b-b28a-4a9b-aa7e-18df09bf383c.cf05ccea05cfdc281aff4fe5748e88d52fbec87d.de-de.xlf"
) | Out-Null
$ws2.Range("F3").Font.Underline = $true
$ws2.Range("F3").Font.Color = 15570276

$ws2.Range("G3").Value = "2016-01-18 07:09:09"
